# Testes unitários adicionados e comentários
#
# The underlying simulation (task scheduler report) was re-run, producing new
# timing numbers for each task/class row on "Tabela de Tarefas" and new
# aggregate statistics on "Análise de Escalonamento". The row/label text is
# unchanged - only the numeric measurements move.

$wb = $excel.ActiveWorkbook

$wsTarefas = $wb.Worksheets.Item("Tabela de Tarefas")
$wsAnalise = $wb.Worksheets.Item("Análise de Escalonamento")

# --- "Tabela de Tarefas": columns are Tarefas | Classe | Ji | Ci | Pi | Di ---
# Row 2 - Company
$wsTarefas.Range("C2").Value = 8868
$wsTarefas.Range("D2").Value = 129
$wsTarefas.Range("E2").Value = 130

# Row 3 - AlphaBank
$wsTarefas.Range("C3").Value = 6634
$wsTarefas.Range("D3").Value = 2026
$wsTarefas.Range("E3").Value = 2026

# Row 4 - Driver
$wsTarefas.Range("C4").Value = 8885
$wsTarefas.Range("D4").Value = 206
$wsTarefas.Range("E4").Value = 311

# Row 5 - Car
$wsTarefas.Range("C5").Value = 8883
$wsTarefas.Range("D5").Value = 577
$wsTarefas.Range("E5").Value = 687

# Row 6 - CarManipulator
$wsTarefas.Range("C6").Value = 8997
$wsTarefas.Range("D6").Value = 127
$wsTarefas.Range("E6").Value = 128

# Row 7 - CompanyAttExcel
$wsTarefas.Range("C7").Value = 8874
$wsTarefas.Range("D7").Value = 68
$wsTarefas.Range("E7").Value = 757

# Row 8 - AlphaBankAttExcel (Ci unchanged at 18)
$wsTarefas.Range("C8").Value = 6636
$wsTarefas.Range("E8").Value = 2042

# Row 9 - ExecutaSimulador
$wsTarefas.Range("C9").Value = 6632
$wsTarefas.Range("D9").Value = 220
$wsTarefas.Range("E9").Value = 220

# Row 10 - SpendFuel
$wsTarefas.Range("C10").Value = 8994
$wsTarefas.Range("D10").Value = 1217
$wsTarefas.Range("E10").Value = 1218

# Row 11 - Account
$wsTarefas.Range("C11").Value = 8991
$wsTarefas.Range("D11").Value = 2016
$wsTarefas.Range("E11").Value = 2016

# --- "Análise de Escalonamento": aggregate stats from the re-run ---
$wsAnalise.Range("B1").Value = 3259956199
$wsAnalise.Range("B2").Value = 7.584580303991242
$wsAnalise.Range("B3").Value = 9.542250513793233
# B4 ("Escalonável" -> "Sim") is unchanged.
